$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 15. This shifts the old row 15 (totals)
# --- and row 16 (footer) down to 16 and 17, and shifts their merged cells
# --- along with them.
$ws.Rows.Item(15).Insert()

# Row heights: the newly inserted blank row takes the old row 15's height
# (24.75), while the totals row (now 16) grows to 25.5.
$ws.Range("A15:Q15").RowHeight = 24.75
$ws.Range("A16:Q16").RowHeight = 25.5

# Recreate the merged ranges for the new product row (mirrors rows 7-14).
$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

# Copy the formatting (fonts, fills, borders, alignment) of an existing
# product row into the new row 15 so it matches the table's styling.
$ws.Range("A7:Q7").Copy($ws.Range("A15:Q15"))

# --- Populate the new product row (item #9: كالونا) ---
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "كالونا "
$ws.Range("H15").Value = "0:0"

$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "0"
$ws.Range("L15").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N15").Value = "15.00"

$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "15.0000"
$ws.Range("P15").NumberFormat = "0.00"

$ws.Range("Q15").Value = "1:0"

# --- Update the totals row (now row 16) with the new grand total ---
$ws.Range("P16").Value = 439.91000000000002

# --- Update the generated-timestamp footer (now row 17) ---
$ws.Range("A17").Value = "Wednesday, 10 September, 2025 10:06 AM"
